# Update countries & provincias Spain
#
# This script applies the COVID data refresh represented by the diff:
#  1. Reorders "Estado de Palestina" and "Camboya" up in the countries
#     list (right after Uzbekistan, before Mauricio) and updates the
#     case numbers for the affected rows so every country keeps its own
#     correct statistics after the shuffle.
#  2. Bumps the "Datos actualizados" timestamp from 03:50 to 04:20.
#  3. Refreshes case counts for Estados Unidos, Corea del Sur and Noruega.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "last updated" timestamp (row 1) ---------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 04:20"

# --- 2. Refresh existing country rows with new case counts -----------------
# Estados Unidos (row 4)
$ws.Range("B4").Value = 123750
$ws.Range("C4").Value = 172
$ws.Range("E4").Value = 118292
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 2227

# Corea del Sur (row 14)
$ws.Range("B14").Value = 9583
$ws.Range("C14").Value = 105
$ws.Range("D14").Value = 5033
$ws.Range("E14").Value = 4398
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = 152

# Noruega (row 20)
$ws.Range("B20").Value = 4032
$ws.Range("C20").Value = 17
$ws.Range("E20").Value = 4002

# --- 3. Reorder Estado de Palestina & Camboya ahead of Mauricio ------------
# Rows 104-108 previously held (in order): Mauricio, Guadalupe,
# Costa de Marfil, Camboya, Estado de Palestina. They now need to read
# (in order): Estado de Palestina, Camboya, Mauricio, Guadalupe,
# Costa de Marfil - each carrying its own refreshed statistics.

# Row 104: now Estado de Palestina
$ws.Range("A104").Value = "Estado de Palestina"
$ws.Range("B104").Value = 104
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 18
$ws.Range("E104").Value = 85
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 1

# Row 105: now Camboya
$ws.Range("A105").Value = "Camboya"
$ws.Range("B105").Value = 103
$ws.Range("C105").Value = 4
$ws.Range("D105").Value = 21
$ws.Range("E105").Value = 82
$ws.Range("F105").Value = 1
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 0

# Row 106: now Mauricio
$ws.Range("A106").Value = "Mauricio"
$ws.Range("B106").Value = 102
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 100
$ws.Range("F106").Value = 1
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 2

# Row 107: now Guadalupe
$ws.Range("A107").Value = "Guadalupe"
$ws.Range("B107").Value = 102
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 17
$ws.Range("E107").Value = 83
$ws.Range("F107").Value = 4
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 2

# Row 108: now Costa de Marfil
$ws.Range("A108").Value = "Costa de Marfil"
$ws.Range("B108").Value = 101
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 3
$ws.Range("E108").Value = 98
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 0
